# Updated cryptos list (Price / Volume(1h) columns) to match the latest
# scraped snapshot. Values are written as plain text (not numbers) so that
# formats such as "41.787.59" (thousand-grouped) or "1.00" (trailing zero)
# round-trip exactly, matching how the source data is stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.787.59"
Set-TextValue $ws.Range("E2") "  +0.63%  "
Set-TextValue $ws.Range("D3") "2.477.26"
Set-TextValue $ws.Range("E3") "  +0.57%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "319.13"
Set-TextValue $ws.Range("E5") "  +1.59%  "
Set-TextValue $ws.Range("D6") "93.24"
Set-TextValue $ws.Range("E6") "  +2.12%  "
Set-TextValue $ws.Range("E7") "  +0.81%  "
Set-TextValue $ws.Range("E8") "  +0.03%  "
Set-TextValue $ws.Range("E9") "  +1.25%  "
Set-TextValue $ws.Range("D10") "0.0884"
Set-TextValue $ws.Range("E10") "  +11.53%  "
Set-TextValue $ws.Range("D11") "33.28"
Set-TextValue $ws.Range("E11") "  +2.45%  "
Set-TextValue $ws.Range("E12") "  +0.61%  "
Set-TextValue $ws.Range("D13") "2.859.40"
Set-TextValue $ws.Range("E13") "  +0.48%  "
Set-TextValue $ws.Range("D14") "6.92"
Set-TextValue $ws.Range("E14") "  +1.26%  "
Set-TextValue $ws.Range("D15") "15.68"
Set-TextValue $ws.Range("E15") "  -0.59%  "
Set-TextValue $ws.Range("D16") "2.466.92"
Set-TextValue $ws.Range("E16") "  +0.95%  "
Set-TextValue $ws.Range("D17") "0.801"
Set-TextValue $ws.Range("E17") "  +3.53%  "
Set-TextValue $ws.Range("D18") "41.741.64"
Set-TextValue $ws.Range("E18") "  +0.51%  "
Set-TextValue $ws.Range("E19") "  +1.61%  "
Set-TextValue $ws.Range("E20") "  -0.10%  "
Set-TextValue $ws.Range("D21") "71.25"
Set-TextValue $ws.Range("E21") "  +0.49%  "
Set-TextValue $ws.Range("E22") "  +2.51%  "
Set-TextValue $ws.Range("D23") "242.21"
Set-TextValue $ws.Range("E23") "  +1.91%  "
Set-TextValue $ws.Range("D24") "2.76"
Set-TextValue $ws.Range("E24") "  +1.81%  "
Set-TextValue $ws.Range("E25") "  +2.70%  "
Set-TextValue $ws.Range("E26") "  +0.00%  "
Set-TextValue $ws.Range("D27") "25.28"
Set-TextValue $ws.Range("E27") "  +3.30%  "
Set-TextValue $ws.Range("D28") "2.26"
Set-TextValue $ws.Range("E28") "  +0.94%  "
Set-TextValue $ws.Range("D29") "9.76"
Set-TextValue $ws.Range("E29") "  +1.07%  "
Set-TextValue $ws.Range("D30") "36.95"
Set-TextValue $ws.Range("E30") "  +4.76%  "
Set-TextValue $ws.Range("D31") "158.65"
Set-TextValue $ws.Range("E31") "  +1.63%  "
Set-TextValue $ws.Range("E32") "  +1.98%  "
Set-TextValue $ws.Range("E33") "  -0.12%  "
Set-TextValue $ws.Range("D34") "0.0766"
Set-TextValue $ws.Range("E34") "  +1.06%  "
Set-TextValue $ws.Range("D35") "2.56"
Set-TextValue $ws.Range("E35") "  -0.31%  "
Set-TextValue $ws.Range("D36") "17.43"
Set-TextValue $ws.Range("E36") "  +1.23%  "
Set-TextValue $ws.Range("E37") "  +5.49%  "
Set-TextValue $ws.Range("E38") "  +1.96%  "
Set-TextValue $ws.Range("E39") "  +1.83%  "
Set-TextValue $ws.Range("D40") "0.105"
Set-TextValue $ws.Range("E40") "  +1.66%  "
Set-TextValue $ws.Range("E41") "  +1.37%  "
Set-TextValue $ws.Range("E42") "  +8.12%  "
Set-TextValue $ws.Range("D43") "2.007.91"
Set-TextValue $ws.Range("E43") "  +3.36%  "
Set-TextValue $ws.Range("D44") "19.38"
Set-TextValue $ws.Range("E44") "  +3.96%  "
Set-TextValue $ws.Range("E45") "  +1.07%  "
Set-TextValue $ws.Range("D46") "2.98"
Set-TextValue $ws.Range("E46") "  +3.16%  "
Set-TextValue $ws.Range("D47") "9.49"
Set-TextValue $ws.Range("E47") "  +5.40%  "
Set-TextValue $ws.Range("D48") "2.716.29"
Set-TextValue $ws.Range("E48") "  +0.40%  "
Set-TextValue $ws.Range("D49") "77.07"
Set-TextValue $ws.Range("E49") "  +8.20%  "
Set-TextValue $ws.Range("D50") "98.05"
Set-TextValue $ws.Range("E50") "  +1.52%  "
Set-TextValue $ws.Range("D51") "67.41"
Set-TextValue $ws.Range("E51") "  +1.04%  "
